$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ACCOUNT sheet: add organization accounts (ACC_CREATE_ORG) ---
# Row 6: existing phone-style account becomes organization code 02001
$ws.Range("A6").Value = "'02001"
$ws.Range("B6").Value = "'02001"
$ws.Range("C6").Value = "1"
$ws.Range("D6").Value = "1"

# Row 7: organization code 02002
$ws.Range("A7").Value = "02002"
$ws.Range("B7").Value = "02002"
$ws.Range("C7").Value = "1"
$ws.Range("D7").Value = "1"

# Row 8: organization code 02003
$ws.Range("A8").Value = "'02003"
$ws.Range("B8").Value = "'02003"
$ws.Range("C8").Value = "1"
$ws.Range("D8").Value = "1"

# Row 9 (new): organization code 02004
$ws.Range("A9").Value = "02004"
$ws.Range("B9").Value = "02004"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "1"

# Row 10 (new): organization code 02005
$ws.Range("A10").Value = "'02005"
$ws.Range("B10").Value = "'02005"
$ws.Range("C10").Value = "1"
$ws.Range("D10").Value = "1"

# Update the active selection to reflect where the edit was last made
$ws.Activate()
$excel.Application.Goto($ws.Range("D6:D10"))
